# Sprint Backlog update — fill in the previously-blank "second task" rows
# for each user story with the assignee name and the Week1/Week2/Remaining
# numbers, matching the rest of the backlog table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# User Story: "User - Create a Trip ..." -> second task row (row 4) is
# worked by the same assignee as row 3 (Jamia).
$ws.Range("B4").Value = "Jamia"
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0

# User Story: "User - Add Waypoints ..." -> second task row (row 6) is
# worked by the same assignee as rows 5/7 (Austin).
$ws.Range("B6").Value = "Austin"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0

# User Story: "User - I want to log into the system" -> second task row
# (row 10) is worked by the same assignee as row 9 (Brian).
$ws.Range("B10").Value = "Brian"
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0
